$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groupes.slk")

# Update the "IdxPar" column values (A2:A11) from 9040000000 to 2010000000
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = 2010000000
}

# Update the selected cell / active cell in the sheet view from D1 to E11
$ws.Range("E11").Select()

# Turn off concurrent calculation (calcPr concurrentCalc="0")
$excel.EnableConcurrentCalculation = $false
